$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.885.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "'2.445.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'578.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'141.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.85%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'2.440.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").Value = "'25.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").Value = "'2.883.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "'61.832.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "'2.437.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'10.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.28%  "
$ws.Range("D20").Value = "'7.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "'324.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'4.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'1.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.87%  "
$ws.Range("D25").Value = "'64.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").Value = "'9.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").Value = "'584.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.72%  "
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'0.0₃0931"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").Value = "'7.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").Value = "'1.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'4.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.01%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.372"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'152.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("D40").Value = "'18.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").Value = "'5.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.77%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'42.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  -4.45%  "
$ws.Range("E45").Value = "  -5.06%  "
$ws.Range("D46").Value = "'0.0₆0279"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +18.00%  "
$ws.Range("D47").Value = "'141.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").Value = "'3.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("D51").Value = "'19.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.82%  "
